# Update gh-pages output data (generated at 456a3b4)
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - F column updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6732
$ws1.Range("F4").Value = 110
$ws1.Range("F5").Value = 157
$ws1.Range("F7").Value = 81
$ws1.Range("F8").Value = 590
$ws1.Range("F9").Value = 47

# Sheet "全部类型" (All types) - F and G column updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6732
$ws4.Range("G2").Value = 58
$ws4.Range("F5").Value = 110
$ws4.Range("F6").Value = 157
$ws4.Range("F9").Value = 81
$ws4.Range("F10").Value = 590
$ws4.Range("F11").Value = 47

$wb.Save()
